# Atualizado por script em 20-11-2023 08:45
#
# This script reorders several match rows (the F:V "match detail" columns —
# home/away teams, scores, odds, timestamps and url — while columns A:E,
# the index / country / tournament / season / kickoff date, stay put) and
# appends 5 new match rows (119-123) at the bottom of the sheet, matching
# the source scrape's freshly-pulled rows for round 2023-11-18.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $rowA, $rowB) {
    $valsA = @()
    $valsB = @()
    for ($c = 6; $c -le 22; $c++) {
        $valsA += ,$ws.Cells.Item($rowA, $c).Value()
        $valsB += ,$ws.Cells.Item($rowB, $c).Value()
    }
    for ($i = 0; $i -lt $valsA.Length; $i++) {
        $c = 6 + $i
        $ws.Cells.Item($rowA, $c).Value = $valsB[$i]
        $ws.Cells.Item($rowB, $c).Value = $valsA[$i]
    }
}

# Simple pairwise swaps of the match-detail columns (F:V) between two rows.
Swap-Rows $ws 62 64
Swap-Rows $ws 77 78
Swap-Rows $ws 91 92
Swap-Rows $ws 103 105

# Rows 108/109/110 rotate: new108 = old110, new109 = old108, new110 = old109.
$r108 = @()
$r109 = @()
$r110 = @()
for ($c = 6; $c -le 22; $c++) {
    $r108 += ,$ws.Cells.Item(108, $c).Value()
    $r109 += ,$ws.Cells.Item(109, $c).Value()
    $r110 += ,$ws.Cells.Item(110, $c).Value()
}
for ($i = 0; $i -lt $r108.Length; $i++) {
    $c = 6 + $i
    $ws.Cells.Item(108, $c).Value = $r110[$i]
    $ws.Cells.Item(109, $c).Value = $r108[$i]
    $ws.Cells.Item(110, $c).Value = $r109[$i]
}

# Append 5 freshly-scraped rows (118-122 in Indice terms -> sheet rows 119-123).
$newRows = @(
    @{ Idx=118; Date=45248.54166666666; F="Sloboda";      G=0; H="Tekstilac Odzaci";       I=2;
       J=2.67; K="18/11/2023 02:12"; L=2.81; M="18/11/2023 12:55";
       N=2.63; O="18/11/2023 02:12"; P=2.77; Q="18/11/2023 12:55";
       R=2.73; S="18/11/2023 02:12"; T=2.57; U="18/11/2023 12:55";
       V="https://www.betexplorer.com/football/serbia/prva-liga/sloboda-tekstilac-odzaci/UupG5x4L/" },
    @{ Idx=119; Date=45248.54166666666; F="Jedinstvo U."; G=3; H="Graficar Beograd";       I=1;
       J=1.98; K="18/11/2023 02:12"; L=1.95; M="18/11/2023 12:51";
       N=3.07; O="18/11/2023 02:12"; P=3.27; Q="18/11/2023 12:51";
       R=3.45; S="18/11/2023 02:12"; T=3.49; U="18/11/2023 12:51";
       V="https://www.betexplorer.com/football/serbia/prva-liga/jedinstvo-ub-graficar-beograd/dSo87bZ8/" },
    @{ Idx=120; Date=45248.54166666666; F="Metalac";      G=2; H="OFK Beograd";            I=0;
       J=4.26; K="18/11/2023 02:12"; L=4.21; M="18/11/2023 12:55";
       N=3.19; O="18/11/2023 02:12"; P=3.32; Q="18/11/2023 12:55";
       R=1.74; S="18/11/2023 02:12"; T=1.76; U="18/11/2023 12:55";
       V="https://www.betexplorer.com/football/serbia/prva-liga/metalac-ofk-beograd/jkoC6IlF/" },
    @{ Idx=121; Date=45248.54166666666; F="Vrsac";        G=1; H="Radnicki S. Mitrovica";  I=0;
       J=2.21; K="18/11/2023 02:12"; L=2.43; M="18/11/2023 12:30";
       N=2.78; O="18/11/2023 02:12"; P=2.64; Q="18/11/2023 12:30";
       R=3.25; S="18/11/2023 02:12"; T=3.18; U="18/11/2023 12:30";
       V="https://www.betexplorer.com/football/serbia/prva-liga/vrsac-radnicki-s-mitrovica/Ai1D4dJR/" },
    @{ Idx=122; Date=45248.625;         F="Macva";        G=1; H="Kolubara";               I=0;
       J=2.11; K="18/11/2023 03:13"; L=2.25; M="18/11/2023 14:58";
       N=2.85; O="18/11/2023 03:13"; P=2.66; Q="18/11/2023 14:58";
       R=3.37; S="18/11/2023 03:13"; T=3.52; U="18/11/2023 14:58";
       V="https://www.betexplorer.com/football/serbia/prva-liga/macva-sabac-kolubara/x6RZQdn5/" }
)

$destRow = 119
foreach ($row in $newRows) {
    # Clone formatting (bold/bordered index style, datetime number format)
    # from the previous row so the appended rows look identical to the rest
    # of the table instead of picking up default General formatting.
    $ws.Range("A" + ($destRow - 1)).Copy()
    $ws.Range("A" + $destRow).PasteSpecial(-4122)
    $ws.Range("E" + ($destRow - 1)).Copy()
    $ws.Range("E" + $destRow).PasteSpecial(-4122)

    $ws.Cells.Item($destRow, 1).Value = $row.Idx
    $ws.Cells.Item($destRow, 2).Value = "serbia"
    $ws.Cells.Item($destRow, 3).Value = "prva-liga"
    $ws.Cells.Item($destRow, 4).Value = "2023-2024"
    $ws.Cells.Item($destRow, 5).Value = $row.Date
    $ws.Cells.Item($destRow, 6).Value = $row.F
    $ws.Cells.Item($destRow, 7).Value = $row.G
    $ws.Cells.Item($destRow, 8).Value = $row.H
    $ws.Cells.Item($destRow, 9).Value = $row.I
    $ws.Cells.Item($destRow, 10).Value = $row.J
    $ws.Cells.Item($destRow, 11).Value = $row.K
    $ws.Cells.Item($destRow, 12).Value = $row.L
    $ws.Cells.Item($destRow, 13).Value = $row.M
    $ws.Cells.Item($destRow, 14).Value = $row.N
    $ws.Cells.Item($destRow, 15).Value = $row.O
    $ws.Cells.Item($destRow, 16).Value = $row.P
    $ws.Cells.Item($destRow, 17).Value = $row.Q
    $ws.Cells.Item($destRow, 18).Value = $row.R
    $ws.Cells.Item($destRow, 19).Value = $row.S
    $ws.Cells.Item($destRow, 20).Value = $row.T
    $ws.Cells.Item($destRow, 21).Value = $row.U
    $ws.Cells.Item($destRow, 22).Value = $row.V

    $destRow += 1
}
